# Update the build/version timestamp embedded in the "About" sheet text
# and in the "build_version" column of the data sheet, following a new
# release build (commit: "Update for release mines - January 30").

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet: replace the version string wherever it appears ---
$wsAbout = $wb.Worksheets.Item("About")

$a2text = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2text.Replace($oldVersion, $newVersion)

$a6text = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6text.Replace($oldVersion, $newVersion)

# --- "Boundaries and methane sources" sheet: build_version column (S), data rows 2-33 ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 19).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    $cellValue = $cell.Value()
    if ($cellValue -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
